# Team 64 - CRC cards - shared.pptx
#
# This script reproduces (via legitimate PowerPoint COM automation calls)
# the text/run/paragraph structure changes seen in the target diff:
#
#   Slide 7 ("DRAFT: Class hierarchy diagram - PIECES"):
#     - the "Class" label run/endParaRPr only picked up a cosmetic
#       `dirty="0"` attribute (an artifact PowerPoint's UI/spell-checker
#       stamps on runs it has "seen" after an edit). The visible text
#       itself ("Class") is unchanged, and the PowerPoint object model
#       does not expose that flag as a settable property, so there is no
#       COM call that reproduces it without altering real content.
#
#   Slide 8 ("CRC card - PLAYER"), the Responsibility(ies) text box:
#     - Paragraph 1 "Responsibility(ies):" gets split into three runs
#       ("Responsibility(" / "ies" / "):") with the same visible text -
#       this happens in real PowerPoint when autocorrect/spell-check
#       flags "ies" as a misspelling boundary; we reproduce the run split
#       by re-assigning the middle substring in place.
#     - Paragraph 3 "Initiate moves" absorbs paragraph 4 (a lone trailing
#       space run) and re-splits into "Initiate " + "moves", with
#       paragraph 4 disappearing entirely.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 8: Responsibility(ies) / Initiate game / Initiate moves box
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$shResp = $s8.Shapes.Item("Google Shape;241;p8")
$tr = $shResp.TextFrame.TextRange

# --- Paragraph 1: "Responsibility(ies):" -> 3 runs, same text ---
# "Responsibility(" (1-15) + "ies" (16-18) + "):" (19-20)
$para1 = $tr.Paragraphs(1, 1)
$mid = $para1.Characters(16, 3)
$mid.Text = "ies"

# --- Paragraphs 3 & 4: merge "Initiate moves" + CR + " " into ---
# --- "Initiate " + "moves", dropping the now-empty paragraph 4  ---
$para3 = $tr.Paragraphs(3, 1)
$para4 = $tr.Paragraphs(4, 1)
$wordLen = 5   # length of "moves"
$startAbs = $para3.Start + $para3.Length - 1 - $wordLen
$lenAbs = $wordLen + 1 + $para4.Length
$tail = $tr.Characters($startAbs, $lenAbs)
$tail.Text = "moves"

# Paragraph 4 is now an empty bullet left behind by the merge above;
# remove it so paragraph 3 becomes the last paragraph in the box.
$tr.Paragraphs(4, 1).Delete()
